{"js": "// Replace the date line and every \"A\u00d7B=C\" answer cell with its updated value.\n// Every \"old\" string below is unique within the document, so a literal,\n// case-sensitive search-and-replace (one result expected per search) is\n// sufficient and avoids any cross-substitution even though some new values\n// happen to equal other old values elsewhere in the document.\nconst replacements = [\n  [\"2026-02-15 Sunday\", \"2026-02-16 Monday\"],\n  [\"417\u00d72=834\", \"822\u00d74=3288\"],\n  [\"446\u00d73=1338\", \"226\u00d72=452\"],\n  [\"187\u00d75=935\", \"253\u00d76=1518\"],\n  [\"288\u00d72=576\", \"901\u00d73=2703\"],\n  [\"361\u00d73=1083\", \"964\u00d77=6748\"],\n  [\"438\u00d73=1314\", \"319\u00d75=1595\"],\n  [\"600\u00d72=1200\", \"343\u00d75=1715\"],\n  [\"965\u00d74=3860\", \"256\u00d72=512\"],\n  [\"555\u00d78=4440\", \"379\u00d75=1895\"],\n  [\"832\u00d73=2496\", \"726\u00d77=5082\"],\n  [\"965\u00d76=5790\", \"640\u00d76=3840\"],\n  [\"467\u00d75=2335\", \"384\u00d77=2688\"],\n  [\"128\u00d79=1152\", \"601\u00d73=1803\"],\n  [\"302\u00d75=1510\", \"348\u00d75=1740\"],\n  [\"885\u00d77=6195\", \"965\u00d76=5790\"],\n  [\"200\u00d79=1800\", \"284\u00d74=1136\"],\n  [\"322\u00d76=1932\", \"823\u00d74=3292\"],\n  [\"104\u00d72=208\", \"479\u00d73=1437\"],\n  [\"777\u00d77=5439\", \"524\u00d73=1572\"],\n  [\"216\u00d76=1296\", \"971\u00d77=6797\"],\n  [\"757\u00d74=3028\", \"261\u00d75=1305\"],\n  [\"651\u00d79=5859\", \"431\u00d75=2155\"],\n  [\"733\u00d75=3665\", \"846\u00d72=1692\"],\n  [\"565\u00d74=2260\", \"366\u00d73=1098\"],\n  [\"858\u00d78=6864\", \"150\u00d78=1200\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Replace only the first (and, given uniqueness, only) match so a later\n  // replacement that happens to introduce an earlier \"old\" value elsewhere\n  // in the list is never re-touched.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"A\u00d7B=C\" answer cell with its updated\n# value. Every \"old\" string is unique within the document, so a literal,\n# case-sensitive Find/Replace (wdReplaceOne = 1) that matches exactly one\n# occurrence is safe -- even though a couple of the \"new\" values happen to\n# equal \"old\" values used elsewhere later in the list, because each\n# replacement only ever touches the single remaining occurrence of its own\n# (still-unique-at-that-point) old text.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2026-02-15 Sunday\", \"2026-02-16 Monday\"),\n    @(\"417\u00d72=834\", \"822\u00d74=3288\"),\n    @(\"446\u00d73=1338\", \"226\u00d72=452\"),\n    @(\"187\u00d75=935\", \"253\u00d76=1518\"),\n    @(\"288\u00d72=576\", \"901\u00d73=2703\"),\n    @(\"361\u00d73=1083\", \"964\u00d77=6748\"),\n    @(\"438\u00d73=1314\", \"319\u00d75=1595\"),\n    @(\"600\u00d72=1200\", \"343\u00d75=1715\"),\n    @(\"965\u00d74=3860\", \"256\u00d72=512\"),\n    @(\"555\u00d78=4440\", \"379\u00d75=1895\"),\n    @(\"832\u00d73=2496\", \"726\u00d77=5082\"),\n    @(\"965\u00d76=5790\", \"640\u00d76=3840\"),\n    @(\"467\u00d75=2335\", \"384\u00d77=2688\"),\n    @(\"128\u00d79=1152\", \"601\u00d73=1803\"),\n    @(\"302\u00d75=1510\", \"348\u00d75=1740\"),\n    @(\"885\u00d77=6195\", \"965\u00d76=5790\"),\n    @(\"200\u00d79=1800\", \"284\u00d74=1136\"),\n    @(\"322\u00d76=1932\", \"823\u00d74=3292\"),\n    @(\"104\u00d72=208\", \"479\u00d73=1437\"),\n    @(\"777\u00d77=5439\", \"524\u00d73=1572\"),\n    @(\"216\u00d76=1296\", \"971\u00d77=6797\"),\n    @(\"757\u00d74=3028\", \"261\u00d75=1305\"),\n    @(\"651\u00d79=5859\", \"431\u00d75=2155\"),\n    @(\"733\u00d75=3665\", \"846\u00d72=1692\"),\n    @(\"565\u00d74=2260\", \"366\u00d73=1098\"),\n    @(\"858\u00d78=6864\", \"150\u00d78=1200\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n\n    # Forward, Wrap=wdFindContinue(1), Format=False, Replace=wdReplaceOne(1)\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n}\n"}
